$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (prices) that must remain plain
# text (matching the workbook's original inlineStr cells) rather than being
# auto-converted to numbers by Excel. We temporarily switch those cells to the
# Text number format, assign the value, then restore the original style so the
# cell formatting is unaffected.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "62.336.74"
$ws.Range("E2").Value = "  +0.76%  "
Set-TextValue $ws.Range("D3") "3.432.43"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws.Range("D5") "413.86"
$ws.Range("E5").Value = "  +1.12%  "
Set-TextValue $ws.Range("D6") "128.83"
$ws.Range("E6").Value = "  +0.27%  "
Set-TextValue $ws.Range("D7") "0.624"
$ws.Range("E7").Value = "  -1.94%  "
Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.66%  "
Set-TextValue $ws.Range("D10") "0.140"
$ws.Range("E10").Value = "  +1.80%  "
Set-TextValue $ws.Range("D11") "42.75"
$ws.Range("E11").Value = "  +0.17%  "
Set-TextValue $ws.Range("D12") "0.0000219"
$ws.Range("E12").Value = "  +9.12%  "
Set-TextValue $ws.Range("D13") "9.20"
$ws.Range("E13").Value = "  +1.49%  "
Set-TextValue $ws.Range("D14") "3.974.09"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  -0.22%  "
Set-TextValue $ws.Range("D16") "20.44"
$ws.Range("E16").Value = "  -3.64%  "
Set-TextValue $ws.Range("D17") "3.422.09"
$ws.Range("E17").Value = "  +0.52%  "
Set-TextValue $ws.Range("D18") "12.71"
$ws.Range("E18").Value = "  +5.28%  "
$ws.Range("E19").Value = "  -0.43%  "
Set-TextValue $ws.Range("D20") "62.379.60"
$ws.Range("E20").Value = "  +0.94%  "
Set-TextValue $ws.Range("D21") "475.62"
$ws.Range("E21").Value = "  +4.84%  "
Set-TextValue $ws.Range("D22") "91.67"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  +3.15%  "
Set-TextValue $ws.Range("D24") "13.11"
$ws.Range("E24").Value = "  +1.71%  "
Set-TextValue $ws.Range("D25") "3.31"
$ws.Range("E25").Value = "  +2.11%  "
Set-TextValue $ws.Range("D26") "9.90"
$ws.Range("E26").Value = "  +13.15%  "
Set-TextValue $ws.Range("D27") "33.33"
$ws.Range("E27").Value = "  -0.51%  "
Set-TextValue $ws.Range("D28") "4.77"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  +2.00%  "
Set-TextValue $ws.Range("D30") "11.85"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  -2.08%  "
Set-TextValue $ws.Range("D34") "40.89"
$ws.Range("E34").Value = "  -4.46%  "
$ws.Range("E35").Value = "  -0.01%  "
Set-TextValue $ws.Range("D36") "57.80"
$ws.Range("E36").Value = "  +8.46%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  +0.04%  "
Set-TextValue $ws.Range("D39") "3.04"
$ws.Range("E39").Value = "  +4.87%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D40") "0.325"
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D41") "0.135"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D43") "144.30"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "2.65"
$ws.Range("E44").Value = "  +9.97%  "
Set-TextValue $ws.Range("D45") "2.06"
$ws.Range("E45").Value = "  +4.38%  "
$ws.Range("E46").Value = "  +2.86%  "
Set-TextValue $ws.Range("D47") "2.43"
$ws.Range("E47").Value = "  +20.95%  "
Set-TextValue $ws.Range("D48") "16.38"
$ws.Range("E48").Value = "  -0.77%  "
Set-TextValue $ws.Range("D49") "0.0₃0536"
$ws.Range("E49").Value = "  +31.43%  "
Set-TextValue $ws.Range("D50") "22.24"
$ws.Range("E50").Value = "  +0.18%  "
Set-TextValue $ws.Range("D51") "112.16"
$ws.Range("E51").Value = "  +5.21%  "
